$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update aerogelThickness value (B3): 20 -> 30
$ws.Range("B3").Value = 30

# Update depth value (B6): 4 -> 6
$ws.Range("B6").Value = 6

# Update mirrorThickness (B10/C10): value 3 -> 0.125, unit mm -> in
$ws.Range("B10").Value = 0.125
$ws.Range("C10").Value = "in"

# Remove the formula in B2, keep it as a plain value of 0.25
$ws.Range("B2").Value = 0.25

# Update the view: scroll back to top-left, select B3 instead of B27
$ws.Range("B3").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
